# Update the "Correspond Handoff Datetime" (col D) and
# "Correspond Handback DateTime" (col G) values on row 5 of the
# zh-cn and de-de report sheets to reflect the new handback run.

$wb = $excel.ActiveWorkbook

$wsZhCn = $wb.Sheets.Item("zh-cn")
$wsZhCn.Range("D5").Value = "2016-01-26 12:21:23"
$wsZhCn.Range("G5").Value = "2016-01-26 12:22:08"

$wsDeDe = $wb.Sheets.Item("de-de")
$wsDeDe.Range("D5").Value = "2016-01-26 12:21:34"
$wsDeDe.Range("G5").Value = "2016-01-26 12:22:29"
